$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "_GoBack" bookmark
# (the "To use Machine Learning ..." objective-answer paragraph) and the
# paragraph that should receive it instead (the "How It Was Solved :"
# heading paragraph). We find them by their visible text rather than by
# a hard-coded index so the script is resilient to minor shifts.

$howSolvedPara = $null
$emptyHeadingPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "How It Was Solved") {
        $howSolvedPara = $p
        # The paragraph immediately before it is the empty Heading-1
        # paragraph that needs to be removed.
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.ParagraphStyle.NameLocal -eq "Heading 1") {
            $emptyHeadingPara = $prev
        }
        break
    }
}

# Move the "_GoBack" bookmark onto the start of the "How It Was Solved :"
# paragraph. Re-adding a bookmark with the same name relocates it (Word
# bookmark names are unique per document), so the copy that currently
# sits on the objective-answer paragraph is implicitly removed.
$target = $howSolvedPara.Range.Duplicate
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)

# Delete the now-redundant empty Heading-1 paragraph that used to sit
# just before "How It Was Solved :".
if ($emptyHeadingPara -ne $null) {
    $emptyHeadingPara.Range.Delete()
}
